$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the header style (bold,
# bordered, centered) from the existing header cell H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-77
$iValues = @(9,9,9,9,9,9,9,7,8,9,9,9,9,7,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,7,7,6,8,8,6,7,7,8,8,7,7,7,9,6,7,8,8,8,9,10,7,7,8,7,7,7,9,8,10,9,9,9,8,8,9,9,4,6,6,5,5,5,4,3,3)
$jValues = @(9,9,9,9,9,9,9,7,9,9,9,9,9,7,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,7,8,7,8,8,6,7,7,8,9,7,7,7,9,6,7,8,8,8,9,10,7,7,8,7,7,8,9,8,10,10,9,9,8,9,9,9,4,6,6,5,5,5,4,3,3)

for ($r = 2; $r -le 77; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
